# Append the "COIN vs. random" run data (rows for sample sizes 6..44) below the
# existing header row. Column A holds the sample-size label (stored as TEXT,
# matching the source data), columns B/C hold the numeric win-rate / average
# game-length figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# data rows: sample size (as text), win rate %, avg. game length
$data = @(
  @("6",  95.40000000000001, 31.548),
  @("8",  95.59999999999999, 31.427),
  @("10", 95.59999999999999, 31.772),
  @("12", 95.09999999999999, 31.441),
  @("14", 95.40000000000001, 31.457),
  @("16", 94,                31.557),
  @("18", 96.2,              31.352),
  @("20", 96.2,              31.678),
  @("22", 95.3,              31.77),
  @("24", 95.8,              31.732),
  @("26", 95.5,              31.338),
  @("28", 94.7,              31.532),
  @("30", 96.5,              31.591),
  @("32", 96.3,              31.199),
  @("34", 95.59999999999999, 31.785),
  @("36", 94.8,              31.35),
  @("38", 95.40000000000001, 31.44),
  @("40", 95.2,              31.396),
  @("42", 95.5,              31.333),
  @("44", 95.2,              31.362)
)

$firstRow = 2
$lastRow = $firstRow + $data.Length - 1

# Give column A the same bold / centered / bordered look already used by the
# B1:C1 header cells, by copying their format onto A2:A21.
$ws.Range("B1").Copy()
$ws.Range("A$firstRow`:A$lastRow").PasteSpecial(-4122)

$r = $firstRow
foreach ($row in $data) {
  # Column A must stay TEXT (the sample-size labels are strings, not
  # numbers). Build the literal string through a throw-away helper cell
  # (ZZ1) via a formula, then paste just the *value* onto the target cell so
  # the style applied above is preserved while the cell becomes text.
  $ws.Range("ZZ1").Formula = '="' + $row[0] + '"'
  $ws.Range("ZZ1").Copy()
  $ws.Range("A$r").PasteSpecial(-4163)

  $ws.Range("B$r").Value = $row[1]
  $ws.Range("C$r").Value = $row[2]

  $r = $r + 1
}

# Clean up the helper cell so it doesn't leak into the saved sheet.
$ws.Range("ZZ1").Clear()

Write-Host "done"
